# Update simulation-result cells B2:F66 with the recorded values
# (commit: "Auslagern von Excel_Auswertung als eigenes Programm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.27
$ws.Range("D2").Value = 1.5
# Row 3
$ws.Range("B3").Value = 1.6
$ws.Range("C3").Value = 0.305
$ws.Range("D3").Value = 1.222222222222222
$ws.Range("E3").Value = 0.025
# Row 4
$ws.Range("B4").Value = 1.25
$ws.Range("C4").Value = 0.38
$ws.Range("D4").Value = 1.105263157894737
$ws.Range("E4").Value = 0.04
# Row 5
$ws.Range("B5").Value = 2.588235294117647
$ws.Range("C5").Value = 0.425
$ws.Range("D5").Value = 1.051546391752577
$ws.Range("E5").Value = 0.08
$ws.Range("F5").Value = 0.005
# Row 6
$ws.Range("B6").Value = 2.227272727272727
$ws.Range("C6").Value = 0.52
$ws.Range("D6").Value = 1.041322314049587
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 0.01
# Row 7
$ws.Range("B7").Value = 2.4375
$ws.Range("C7").Value = 0.58
$ws.Range("D7").Value = 1.034965034965035
$ws.Range("E7").Value = 0.135
$ws.Range("F7").Value = 0.025
# Row 8
$ws.Range("B8").Value = 2.384615384615385
$ws.Range("C8").Value = 0.605
$ws.Range("D8").Value = 1.032258064516129
$ws.Range("E8").Value = 0.16
$ws.Range("F8").Value = 0.035
# Row 9
$ws.Range("B9").Value = 2.511111111111111
$ws.Range("C9").Value = 0.625
$ws.Range("D9").Value = 1.03030303030303
$ws.Range("E9").Value = 0.19
$ws.Range("F9").Value = 0.035
# Row 10
$ws.Range("B10").Value = 2.527272727272727
$ws.Range("C10").Value = 0.62
$ws.Range("D10").Value = 1.022857142857143
$ws.Range("E10").Value = 0.235
$ws.Range("F10").Value = 0.04
# Row 11
$ws.Range("B11").Value = 2.421875
$ws.Range("C11").Value = 0.59
$ws.Range("D11").Value = 1.02247191011236
$ws.Range("E11").Value = 0.28
$ws.Range("F11").Value = 0.04
# Row 12
$ws.Range("B12").Value = 2.397260273972603
$ws.Range("C12").Value = 0.5649999999999999
$ws.Range("D12").Value = 1.021978021978022
$ws.Range("E12").Value = 0.325
$ws.Range("F12").Value = 0.04
# Row 13
$ws.Range("B13").Value = 2.25
$ws.Range("C13").Value = 0.545
$ws.Range("D13").Value = 1.016129032258065
$ws.Range("E13").Value = 0.36
$ws.Range("F13").Value = 0.04
# Row 14
$ws.Range("B14").Value = 2.467391304347826
$ws.Range("C14").Value = 0.49
$ws.Range("D14").Value = 1.01063829787234
$ws.Range("E14").Value = 0.42
$ws.Range("F14").Value = 0.04
# Row 15
$ws.Range("B15").Value = 2.673267326732673
$ws.Range("C15").Value = 0.465
$ws.Range("D15").Value = 1.010416666666667
$ws.Range("E15").Value = 0.465
$ws.Range("F15").Value = 0.04
# Row 16
$ws.Range("B16").Value = 2.663716814159292
$ws.Range("C16").Value = 0.405
$ws.Range("D16").Value = 1.010416666666667
$ws.Range("E16").Value = 0.52
$ws.Range("F16").Value = 0.045
# Row 17
$ws.Range("B17").Value = 2.758064516129032
$ws.Range("C17").Value = 0.35
$ws.Range("D17").Value = 1.010416666666667
$ws.Range("E17").Value = 0.575
$ws.Range("F17").Value = 0.045
# Row 18
$ws.Range("B18").Value = 2.772727272727273
$ws.Range("C18").Value = 0.32
$ws.Range("D18").Value = 1.010309278350515
$ws.Range("E18").Value = 0.615
$ws.Range("F18").Value = 0.045
# Row 19
$ws.Range("B19").Value = 2.685714285714285
$ws.Range("C19").Value = 0.29
$ws.Range("D19").Value = 1.010204081632653
$ws.Range("E19").Value = 0.65
$ws.Range("F19").Value = 0.05
# Row 20
$ws.Range("B20").Value = 2.701388888888889
$ws.Range("C20").Value = 0.27
$ws.Range("D20").Value = 1.010204081632653
$ws.Range("E20").Value = 0.67
$ws.Range("F20").Value = 0.05
# Row 21
$ws.Range("B21").Value = 2.662251655629139
$ws.Range("C21").Value = 0.23
$ws.Range("D21").Value = 1.005102040816326
$ws.Range("E21").Value = 0.7
$ws.Range("F21").Value = 0.055
# Row 22
$ws.Range("B22").Value = 2.651315789473684
$ws.Range("C22").Value = 0.235
$ws.Range("D22").Value = 1.005050505050505
$ws.Range("E22").Value = 0.705
$ws.Range("F22").Value = 0.055
# Row 23
$ws.Range("B23").Value = 2.621794871794872
$ws.Range("C23").Value = 0.215
$ws.Range("D23").Value = 1.005050505050505
$ws.Range("E23").Value = 0.725
$ws.Range("F23").Value = 0.055
# Row 24
$ws.Range("B24").Value = 2.611111111111111
$ws.Range("C24").Value = 0.185
$ws.Range("D24").Value = 1.005050505050505
$ws.Range("E24").Value = 0.755
$ws.Range("F24").Value = 0.055
# Row 25
$ws.Range("B25").Value = 2.6
$ws.Range("C25").Value = 0.175
$ws.Range("D25").Value = 1.005025125628141
$ws.Range("E25").Value = 0.77
$ws.Range("F25").Value = 0.055
# Row 26
$ws.Range("B26").Value = 2.590361445783132
$ws.Range("C26").Value = 0.17
$ws.Range("D26").Value = 1.005025125628141
$ws.Range("E26").Value = 0.775
$ws.Range("F26").Value = 0.055
# Row 27
$ws.Range("B27").Value = 2.558823529411764
$ws.Range("C27").Value = 0.155
$ws.Range("D27").Value = 1.005
$ws.Range("E27").Value = 0.79
$ws.Range("F27").Value = 0.06
# Row 28
$ws.Range("B28").Value = 2.558823529411764
$ws.Range("C28").Value = 0.155
$ws.Range("D28").Value = 1.005
$ws.Range("E28").Value = 0.79
$ws.Range("F28").Value = 0.06
# Row 29
$ws.Range("B29").Value = 2.540229885057471
$ws.Range("C29").Value = 0.135
$ws.Range("D29").Value = 1.005
$ws.Range("E29").Value = 0.805
$ws.Range("F29").Value = 0.065
# Row 30
$ws.Range("B30").Value = 2.51685393258427
$ws.Range("C30").Value = 0.115
$ws.Range("D30").Value = 1.005
$ws.Range("E30").Value = 0.825
$ws.Range("F30").Value = 0.065
# Row 31
$ws.Range("B31").Value = 2.540983606557377
$ws.Range("C31").Value = 0.08500000000000001
$ws.Range("D31").Value = 1
$ws.Range("E31").Value = 0.85
$ws.Range("F31").Value = 0.065
# Row 32
$ws.Range("B32").Value = 2.527173913043478
$ws.Range("C32").Value = 0.08
$ws.Range("D32").Value = 1
$ws.Range("E32").Value = 0.855
$ws.Range("F32").Value = 0.065
# Row 33
$ws.Range("B33").Value = 2.505376344086022
$ws.Range("C33").Value = 0.07000000000000001
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 0.865
$ws.Range("F33").Value = 0.065
# Row 34
$ws.Range("B34").Value = 2.505376344086022
$ws.Range("C34").Value = 0.07000000000000001
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0.865
$ws.Range("F34").Value = 0.065
# Row 35
$ws.Range("B35").Value = 2.505376344086022
$ws.Range("C35").Value = 0.07000000000000001
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0.865
$ws.Range("F35").Value = 0.065
# Row 36
$ws.Range("B36").Value = 2.505376344086022
$ws.Range("C36").Value = 0.07000000000000001
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 0.865
$ws.Range("F36").Value = 0.065
# Row 37
$ws.Range("B37").Value = 2.492063492063492
$ws.Range("C37").Value = 0.055
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0.88
$ws.Range("F37").Value = 0.065
# Row 38
$ws.Range("B38").Value = 2.492063492063492
$ws.Range("C38").Value = 0.055
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 0.88
$ws.Range("F38").Value = 0.065
# Row 39
$ws.Range("B39").Value = 2.492063492063492
$ws.Range("C39").Value = 0.055
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 0.88
$ws.Range("F39").Value = 0.065
# Row 40
$ws.Range("B40").Value = 2.489473684210526
$ws.Range("C40").Value = 0.05
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 0.885
$ws.Range("F40").Value = 0.065
# Row 41
$ws.Range("B41").Value = 2.505208333333333
$ws.Range("C41").Value = 0.04
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 0.895
$ws.Range("F41").Value = 0.065
# Row 42
$ws.Range("B42").Value = 2.556701030927835
$ws.Range("C42").Value = 0.03
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = 0.905
$ws.Range("F42").Value = 0.065
# Row 43
$ws.Range("B43").Value = 2.548717948717949
$ws.Range("C43").Value = 0.025
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 0.91
$ws.Range("F43").Value = 0.065
# Row 44
$ws.Range("B44").Value = 2.548717948717949
$ws.Range("C44").Value = 0.025
$ws.Range("D44").Value = 1
$ws.Range("E44").Value = 0.91
$ws.Range("F44").Value = 0.065
# Row 45
$ws.Range("B45").Value = 2.548717948717949
$ws.Range("C45").Value = 0.025
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 0.91
$ws.Range("F45").Value = 0.065
# Row 46
$ws.Range("B46").Value = 2.548717948717949
$ws.Range("C46").Value = 0.025
$ws.Range("D46").Value = 1
$ws.Range("E46").Value = 0.91
$ws.Range("F46").Value = 0.065
# Row 47
$ws.Range("B47").Value = 2.535714285714286
$ws.Range("C47").Value = 0.02
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = 0.915
$ws.Range("F47").Value = 0.065
# Row 48
$ws.Range("B48").Value = 2.535714285714286
$ws.Range("C48").Value = 0.02
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = 0.915
$ws.Range("F48").Value = 0.065
# Row 49
$ws.Range("B49").Value = 2.535714285714286
$ws.Range("C49").Value = 0.02
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = 0.915
$ws.Range("F49").Value = 0.065
# Row 50
$ws.Range("B50").Value = 2.515151515151515
$ws.Range("C50").Value = 0.01
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 0.925
$ws.Range("F50").Value = 0.065
# Row 51
$ws.Range("B51").Value = 2.515151515151515
$ws.Range("C51").Value = 0.01
$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 0.925
$ws.Range("F51").Value = 0.065
# Row 52
$ws.Range("B52").Value = 2.515151515151515
$ws.Range("C52").Value = 0.01
$ws.Range("D52").Value = 1
$ws.Range("E52").Value = 0.925
$ws.Range("F52").Value = 0.065
# Row 53
$ws.Range("B53").Value = 2.515151515151515
$ws.Range("C53").Value = 0.01
$ws.Range("D53").Value = 1
$ws.Range("E53").Value = 0.925
$ws.Range("F53").Value = 0.065
# Row 54
$ws.Range("B54").Value = 2.50251256281407
$ws.Range("C54").Value = 0.005
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 0.93
$ws.Range("F54").Value = 0.065
# Row 55
$ws.Range("B55").Value = 2.50251256281407
$ws.Range("C55").Value = 0.005
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 0.93
$ws.Range("F55").Value = 0.065
# Row 56
$ws.Range("B56").Value = 2.50251256281407
$ws.Range("C56").Value = 0.005
$ws.Range("D56").Value = 1
$ws.Range("E56").Value = 0.93
$ws.Range("F56").Value = 0.065
# Row 57
$ws.Range("B57").Value = 2.50251256281407
$ws.Range("C57").Value = 0.005
$ws.Range("D57").Value = 1
$ws.Range("E57").Value = 0.93
$ws.Range("F57").Value = 0.065
# Row 58
$ws.Range("B58").Value = 2.50251256281407
$ws.Range("C58").Value = 0.005
$ws.Range("D58").Value = 1
$ws.Range("E58").Value = 0.93
$ws.Range("F58").Value = 0.065
# Row 59
$ws.Range("B59").Value = 2.50251256281407
$ws.Range("C59").Value = 0.005
$ws.Range("D59").Value = 1
$ws.Range("E59").Value = 0.93
$ws.Range("F59").Value = 0.065
# Row 60
$ws.Range("B60").Value = 2.50251256281407
$ws.Range("C60").Value = 0.005
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 0.93
$ws.Range("F60").Value = 0.065
# Row 61
$ws.Range("B61").Value = 2.50251256281407
$ws.Range("C61").Value = 0.005
$ws.Range("D61").Value = 1
$ws.Range("E61").Value = 0.93
$ws.Range("F61").Value = 0.065
# Row 62
$ws.Range("B62").Value = 2.50251256281407
$ws.Range("C62").Value = 0.005
$ws.Range("D62").Value = 1
$ws.Range("E62").Value = 0.93
$ws.Range("F62").Value = 0.065
# Row 63
$ws.Range("B63").Value = 2.50251256281407
$ws.Range("C63").Value = 0.005
$ws.Range("D63").Value = 1
$ws.Range("E63").Value = 0.93
$ws.Range("F63").Value = 0.065
# Row 64
$ws.Range("B64").Value = 2.50251256281407
$ws.Range("C64").Value = 0.005
$ws.Range("D64").Value = 1
$ws.Range("E64").Value = 0.93
$ws.Range("F64").Value = 0.065
# Row 65
$ws.Range("B65").Value = 2.50251256281407
$ws.Range("C65").Value = 0.005
$ws.Range("D65").Value = 1
$ws.Range("E65").Value = 0.93
$ws.Range("F65").Value = 0.065
# Row 66
$ws.Range("B66").Value = 2.50251256281407
$ws.Range("C66").Value = 0.005
$ws.Range("D66").Value = 1
$ws.Range("E66").Value = 0.93
$ws.Range("F66").Value = 0.065

Write-Output "Updated 320 cells on Sheet1 (rows 2-66)"
